$wb = $excel.ActiveWorkbook

# --- "Chart" sheet: roll the daily export forward by one day ---
# Drop the oldest day (2025-10-28) by deleting row 2; this shifts every
# later row up by one, which also advances each date by a day and moves
# each day's "No video indexed" / "Video indexed" counts up into the row
# above (matching the new day's lower backlog counts).
$chart = $wb.Worksheets.Item("Chart")
$chart.Rows(2).Delete()

# The shift leaves the last row (now row 88, previously empty because the
# sheet only had 88 rows) to be filled with the newest day's data, plus
# three more new days appended after it.
$chart.Range("A88:A91").NumberFormat = "@"

$chart.Range("A88").Value = "2026-01-23"
$chart.Range("B88").Value = 0
$chart.Range("C88").Value = 0
$chart.Range("D88").Value = 0

$chart.Range("A89").Value = "2026-01-24"
$chart.Range("B89").Value = 0
$chart.Range("C89").Value = 0
$chart.Range("D89").Value = 0

$chart.Range("A90").Value = "2026-01-25"
$chart.Range("B90").Value = 0
$chart.Range("C90").Value = 0
$chart.Range("D90").Value = ""

$chart.Range("A91").Value = "2026-01-26"
$chart.Range("B91").Value = 0
$chart.Range("C91").Value = 0
$chart.Range("D91").Value = ""

# --- "Table" sheet: validation status for the unindexed video updated ---
$table = $wb.Worksheets.Item("Table")
$table.Range("B2").Value = "N/A"
$table.Range("C2").Value = 0
